$wb = $excel.ActiveWorkbook

# --- Login sheet: add a new row of test data (Admin / admin123 / Dashboard) ---
$ws1 = $wb.Sheets.Item("Login")
$ws1.Activate()

# Carry the formatting of the row above down onto the new row (row 5's
# style) by copying it in, then overwrite with the new row's own values.
$ws1.Range("A5:C5").Copy()
$ws1.Range("A6:C6").Insert(-4121)

$ws1.Range("A6").Value = "Admin"
$ws1.Range("B6").Value = "admin123"
$ws1.Range("C6").Value = "Dashboard"

# Leave the selection where the author left it when the file was saved.
$ws1.Range("C19").Select()

# --- Add_User sheet: give it an explicit (portrait) page setup, like Login has ---
$ws2 = $wb.Sheets.Item("Add_User")
$ws2.PageSetup.Orientation = 1
